$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("N2").Value = 1.02
$ws.Range("O2").Value = 1.22
$ws.Range("P2").Value = 1.29
$ws.Range("Q2").Value = 1.22
$ws.Range("S2").Value = 1.22
# Row 3
$ws.Range("F3").Value = 1.72
$ws.Range("G3").Value = 1.78
$ws.Range("H3").Value = 5
$ws.Range("I3").Value = 5.5
$ws.Range("K3").Value = 4.4
$ws.Range("Q3").Value = 1.72
$ws.Range("R3").Value = 1.49
$ws.Range("T3").Value = 1.73
$ws.Range("U3").Value = 2.2
$ws.Range("V3").Value = 1.22
$ws.Range("W3").Value = 2.28
$ws.Range("Y3").Value = 23
$ws.Range("Z3").Value = 980
$ws.Range("AA3").Value = 160
$ws.Range("AC3").Value = 9.6
$ws.Range("AD3").Value = 21
$ws.Range("AE3").Value = 65
$ws.Range("AH3").Value = 18.5
$ws.Range("AI3").Value = 65
$ws.Range("AJ3").Value = 18.5
$ws.Range("AL3").Value = 980
$ws.Range("AM3").Value = 110
$ws.Range("AN3").Value = 8.800000000000001
$ws.Range("AO3").Value = 85
# Row 4
$ws.Range("J4").Value = 2.82
$ws.Range("L4").Value = 1.56
$ws.Range("W4").Value = 1.19
$ws.Range("Z4").Value = 9
# Row 5
$ws.Range("J5").Value = 4
# Row 6
$ws.Range("H6").Value = 1.84
$ws.Range("J6").Value = 3.2
$ws.Range("K6").Value = 4.3
# Row 8
$ws.Range("F8").Value = 1.96
$ws.Range("I8").Value = 4.4
$ws.Range("M8").Value = 1.05
$ws.Range("N8").Value = 4
$ws.Range("O8").Value = 1.26
$ws.Range("R8").Value = 1.41
$ws.Range("T8").Value = 1.68
$ws.Range("V8").Value = 1.29
$ws.Range("AB8").Value = 13
$ws.Range("AI8").Value = 60
$ws.Range("AK8").Value = 25
$ws.Range("AM8").Value = 100
# Row 9
$ws.Range("I9").Value = 2.2
$ws.Range("Q9").Value = 1.85
$ws.Range("S9").Value = 3.15
# Row 10
$ws.Range("F10").Value = 3.75
$ws.Range("G10").Value = 5.3
$ws.Range("I10").Value = 2.06
$ws.Range("J10").Value = 3.35
$ws.Range("S10").Value = 2.58
$ws.Range("V10").Value = 1.94
$ws.Range("AB10").Value = 980
$ws.Range("AE10").Value = 980
$ws.Range("AG10").Value = 980
$ws.Range("AH10").Value = 980
# Row 11
$ws.Range("S11").Value = 2.12
# Row 12
$ws.Range("K12").Value = 3.7
$ws.Range("AE12").Value = 980
$ws.Range("AJ12").Value = 980
$ws.Range("AL12").Value = 980
# Row 14
$ws.Range("F14").Value = 2.3
$ws.Range("I14").Value = 3.15
$ws.Range("M14").Value = 1.03
$ws.Range("V14").Value = 1.46
# Row 15
$ws.Range("G15").Value = 11.5
$ws.Range("J15").Value = 5.2
$ws.Range("K15").Value = 5.4
$ws.Range("L15").Value = 1.4
$ws.Range("M15").Value = 1.06
$ws.Range("N15").Value = 3.75
$ws.Range("O15").Value = 1.32
$ws.Range("P15").Value = 1.95
$ws.Range("Q15").Value = 1.95
$ws.Range("R15").Value = 1.36
$ws.Range("S15").Value = 3.45
$ws.Range("T15").Value = 2.4
$ws.Range("U15").Value = 1.64
$ws.Range("V15").Value = 3.4
$ws.Range("W15").Value = 1.1
$ws.Range("X15").Value = 16.5
$ws.Range("Y15").Value = 7.2
$ws.Range("Z15").Value = 7.2
$ws.Range("AA15").Value = 11
$ws.Range("AB15").Value = 30
$ws.Range("AC15").Value = 12
$ws.Range("AD15").Value = 11
$ws.Range("AE15").Value = 16.5
$ws.Range("AF15").Value = 110
$ws.Range("AG15").Value = 44
$ws.Range("AH15").Value = 36
$ws.Range("AI15").Value = 50
$ws.Range("AJ15").Value = 540
$ws.Range("AK15").Value = 240
$ws.Range("AL15").Value = 200
$ws.Range("AM15").Value = 270
$ws.Range("AN15").Value = 370
$ws.Range("AO15").Value = 7.2
# Row 16
$ws.Range("H16").Value = 1.72
$ws.Range("I16").Value = 1.83
$ws.Range("K16").Value = 5.1
$ws.Range("L16").Value = 1.22
$ws.Range("O16").Value = 1.16
$ws.Range("S16").Value = 2.14
$ws.Range("T16").Value = 1.54
$ws.Range("V16").Value = 2.2
$ws.Range("W16").Value = 1.22
$ws.Range("AH16").Value = 20
# Row 17
$ws.Range("F17").Value = 2.68
$ws.Range("G17").Value = 2.98
$ws.Range("H17").Value = 2.4
$ws.Range("I17").Value = 2.66
$ws.Range("P17").Value = 2.38
$ws.Range("R17").Value = 1.51
$ws.Range("S17").Value = 2.3
$ws.Range("AE17").Value = 34
$ws.Range("AJ17").Value = 65
# Row 18
$ws.Range("F18").Value = 1.68
$ws.Range("G18").Value = 1.86
$ws.Range("K18").Value = 4
$ws.Range("L18").Value = 1.47
$ws.Range("M18").Value = 1.08
$ws.Range("N18").Value = 2.92
$ws.Range("P18").Value = 1.65
$ws.Range("Q18").Value = 2.22
$ws.Range("R18").Value = 1.24
$ws.Range("S18").Value = 3.8
$ws.Range("U18").Value = 1.75
$ws.Range("V18").Value = 1.16
$ws.Range("W18").Value = 2.16
$ws.Range("X18").Value = 13.5
$ws.Range("Y18").Value = 20
$ws.Range("AB18").Value = 8
$ws.Range("AH18").Value = 32
$ws.Range("AJ18").Value = 23
$ws.Range("AK18").Value = 27
$ws.Range("AN18").Value = 18.5
# Row 19
$ws.Range("G19").Value = 2.62
$ws.Range("H19").Value = 3.1
$ws.Range("I19").Value = 3.15
$ws.Range("N19").Value = 3.5
$ws.Range("P19").Value = 1.85
$ws.Range("U19").Value = 2.1
$ws.Range("X19").Value = 12.5
$ws.Range("Z19").Value = 19.5
# Row 20
$ws.Range("F20").Value = 1.85
$ws.Range("H20").Value = 5.4
$ws.Range("J20").Value = 3.6
$ws.Range("L20").Value = 1.51
$ws.Range("T20").Value = 2.2
$ws.Range("W20").Value = 2.16
# Row 21
$ws.Range("H21").Value = 10.5
$ws.Range("K21").Value = 6.4
$ws.Range("M21").Value = 1.05
$ws.Range("O21").Value = 1.25
$ws.Range("P21").Value = 2.2
$ws.Range("Q21").Value = 1.75
$ws.Range("R21").Value = 1.46
$ws.Range("S21").Value = 2.94
$ws.Range("T21").Value = 2.22
$ws.Range("AA21").Value = 530
$ws.Range("AB21").Value = 8.199999999999999
$ws.Range("AC21").Value = 13.5
$ws.Range("AE21").Value = 250
$ws.Range("AH21").Value = 34
$ws.Range("AI21").Value = 150
$ws.Range("AK21").Value = 15
$ws.Range("AM21").Value = 210
$ws.Range("AN21").Value = 5.9
$ws.Range("AO21").Value = 280

$wb.Save()